$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 06:52"

# Reorder Kirguistan ahead of Bolivia/Albania (rows 99-101), with updated
# per-country statistics, mirroring a refreshed data pull that re-sorted
# the table by total cases.
$ws.Range("A99").Value = "Kirguistan"
$ws.Range("B99").Value = 568
$ws.Range("C99").Value = 14
$ws.Range("D99").Value = 201
$ws.Range("E99").Value = 360
$ws.Range("F99").Value = 5
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 7

$ws.Range("A100").Value = "Bolivia"
$ws.Range("B100").Value = 564
$ws.Range("C100").Value = 44
$ws.Range("D100").Value = 31
$ws.Range("E100").Value = 500
$ws.Range("F100").Value = 3
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 33

$ws.Range("A101").Value = "Albania"
$ws.Range("B101").Value = 562
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 314
$ws.Range("E101").Value = 222
$ws.Range("F101").Value = 5
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 26

# Camboya (row 135) stats refresh
$ws.Range("D135").Value = 107
$ws.Range("E135").Value = 15
